$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: 标志 (flag) - same value for all rows
$ws.Range("A2").Value = "CRA20191227194043"
$ws.Range("A3:A11").Value = "CRA20191227194043"

# Column B: 序号 (sequence code) - unique per row
$ws.Range("B2").Value = "CRA201912271940430001"
$ws.Range("B3").Value = "CRA201912271940430002"
$ws.Range("B4").Value = "CRA201912271940430003"
$ws.Range("B5").Value = "CRA201912271940430004"
$ws.Range("B6").Value = "CRA201912271940430005"
$ws.Range("B7").Value = "CRA201912271940430006"
$ws.Range("B8").Value = "CRA201912271940430007"
$ws.Range("B9").Value = "CRA201912271940430008"
$ws.Range("B10").Value = "CRA201912271940430009"
$ws.Range("B11").Value = "CRA201912271940430010"

# Column C: 题名 (title) - unique per row
$ws.Range("C2").Value = "以色列这次真捅了马蜂窝，179枚导弹连番轰炸，这才是真正苦战"
$ws.Range("C3").Value = "卫星曝印度边境突现神秘基地：印军紧急侦察，却发现屏幕一片模糊"
$ws.Range("C4").Value = "乌克兰又在大甩卖？价格十分尴尬，大国不需要，小国买不起"
$ws.Range("C5").Value = "联合国大会上，中方一票否决西方提案，美当初阻扰马达西奇终于遭报复"
$ws.Range("C6").Value = "后悔拒绝中方索赔，波音CEO被炒鱿鱼，留下最后2句忠告"
$ws.Range("C7").Value = "投票再次通过，特朗普正式签名，大局已定，白宫终于松一口气"
$ws.Range("C8").Value = "美国确认制裁俄欧天然气管道，个别欧企立刻停工，俄提前留了一手"
$ws.Range("C9").Value = "韩国制造不靠谱！挪威27500吨巨舰曝重大缺陷，服役不久就被禁航"
$ws.Range("C10").Value = "印巴突然交火，中方有何评论？外交部回应"
$ws.Range("C11").Value = "20万吨！今年中国海军下水吨位又是世界第一"

# Column I: 发表时间 (publish date) - force text, same value for all rows
$dateRange = $ws.Range("I2:I11")
$dateRange.NumberFormat = "@"
$ws.Range("I2").Value = "2019-12-27"
$ws.Range("I3:I11").Value = "2019-12-27"
$dateRange.ClearFormats()

# Column J: 下载地址 (download url) - unique per row, keep hyperlink style (s=2) intact
$ws.Range("J2").Value = "http://www.sohu.com/a/363188949_120147869"
$ws.Range("J3").Value = "http://www.sohu.com/a/363175748_100145375"
$ws.Range("J4").Value = "http://www.sohu.com/a/363150133_120157852"
$ws.Range("J5").Value = "http://www.sohu.com/a/363162539_637401"
$ws.Range("J6").Value = "http://www.sohu.com/a/363189846_120098002"
$ws.Range("J7").Value = "http://www.sohu.com/a/363138524_100018095"
$ws.Range("J8").Value = "http://www.sohu.com/a/363182702_100143135"
$ws.Range("J9").Value = "http://www.sohu.com/a/363174104_120113110"
$ws.Range("J10").Value = "http://www.sohu.com/a/363148524_162522"
$ws.Range("J11").Value = "http://www.sohu.com/a/363182424_115479"

Write-Host "Edit complete"
